# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.245.42"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.823.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.666"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.401"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07088"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008760"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.247.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.315"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.049.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.934"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.267"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.291"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08897"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7776"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.527"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.917"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01962"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05248"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.264"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.91%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.385"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +20.61%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.901"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1690"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.624"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5039"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.673"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
